$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells hold plain text (prices / percentages formatted as strings),
# not numbers -- force Text format so Excel does not auto-convert
# dotted numbers like "26.873.89" or padded percents into numerics.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "26.873.89"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.13"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "310.81"
$ws.Range("E5").Value = "  -1.10%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.47%  "

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4582"
$ws.Range("E7").Value = "  -0.72%  "

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07143"
$ws.Range("E9").Value = "  -2.47%  "

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8710"
$ws.Range("E10").Value = "  -1.27%  "

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07769"
$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "19.51"
$ws.Range("E12").Value = "  -2.12%  "

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "1.821.24"
$ws.Range("E13").Value = "  -2.06%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "6.373"
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "86.84"
$ws.Range("E16").Value = "  -5.47%  "

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008705"
$ws.Range("E18").Value = "  -3.89%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "26.903.71"
$ws.Range("E20").Value = "  -1.70%  "

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "14.43"
$ws.Range("E21").Value = "  -2.50%  "

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "4.981"
$ws.Range("E22").Value = "  -3.00%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "2.008"
$ws.Range("E24").Value = "  +4.15%  "

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "151.20"
$ws.Range("E25").Value = "  -0.62%  "

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "18.17"
$ws.Range("E26").Value = "  -1.00%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.47%  "

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "113.50"
$ws.Range("E28").Value = "  -2.18%  "

$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = "4.914"
$ws.Range("E29").Value = "  -3.78%  "

$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08791"
$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7432"
$ws.Range("E32").Value = "  -3.47%  "

$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D33").Value = "4.470"
$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D34").Value = "1.124"
$ws.Range("E34").Value = "  -4.40%  "

$ws.Range("D35:E35").NumberFormat = "@"
$ws.Range("D35").Value = "2.503"
$ws.Range("E35").Value = "  -5.76%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.86%  "

$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01934"
$ws.Range("E37").Value = "  -1.38%  "

$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05104"
$ws.Range("E38").Value = "  -2.47%  "

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "2.897"
$ws.Range("E39").Value = "  -2.03%  "

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "6.914"
$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4956"
$ws.Range("E41").Value = "  -3.68%  "

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1596"
$ws.Range("E42").Value = "  -2.46%  "

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "8.268"
$ws.Range("E43").Value = "  -1.97%  "

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4668"
$ws.Range("E44").Value = "  -3.21%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.54%  "

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "10.07"
$ws.Range("E46").Value = "  -2.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.57"

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "1.603"
$ws.Range("E48").Value = "  -2.99%  "

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06078"
$ws.Range("E49").Value = "  -2.31%  "

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "64.41"
$ws.Range("E50").Value = "  -2.07%  "

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "36.66"
$ws.Range("E51").Value = "  -0.36%  "

